$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-coerced to a number by
# Excels type inference (single decimal point, digits only) are written with
# an explicit Text number format, then the format is cleared again so the
# cell ends up with no style index -- matching the source file, where none of
# these cells carry an "s" attribute.

$ws.Range('D2').Value = '27.688.14'
$ws.Range('E2').Value = '  +6.02%  '
$ws.Range('D3').Value = '1.735.61'
$ws.Range('E3').Value = '  +4.87%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.79'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +4.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5449'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +3.63%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.004'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2750'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +2.51%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06715'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +5.27%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.95'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +6.68%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07782'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.25%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.698'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.85%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.739.57'
$ws.Range('E13').Value = '  +3.83%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '1.975.40'
$ws.Range('E14').Value = '  +4.91%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5986'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +6.24%  '
$ws.Range('D16').Value = '0.0₅8419'
$ws.Range('E16').Value = '  +2.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '69.35'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +5.58%  '
$ws.Range('D18').Value = '27.730.35'
$ws.Range('E18').Value = '  +6.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '227.53'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +19.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.838'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +2.99%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.004'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.91'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +5.42%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.233'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +4.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.005'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '148.28'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.741'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +14.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1250'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +3.96%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.466'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +2.82%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '17.15'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +7.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05708'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.89%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.315'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +2.97%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.700'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +5.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.519'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +3.97%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.687'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +6.64%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9772'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.93%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.859'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.02%  '
$ws.Range('E37').Value = '  +1.40%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5984'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +3.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01672'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +4.70%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.981'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.16%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8511'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.71%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '1.049.83'
$ws.Range('E42').Value = '  +2.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.004'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.63'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.33%  '
$ws.Range('D45').Value = '1.880.63'
$ws.Range('E45').Value = '  +4.88%  '
$ws.Range('D46').Value = '0.0₈114'
$ws.Range('E46').Value = '  +7.97%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '59.74'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.39%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.312'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +3.49%  '
$ws.Range('E49').Value = '  +2.21%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.003'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05332'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.00%  '
